$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.994.74"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.054.01"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'246.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "'58.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "'0.0780"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'15.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.21%  "
$ws.Range("D13").Value = "'0.882"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("D14").Value = "2.356.30"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'5.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "2.085.96"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "'18.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "36.968.93"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'73.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'5.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'238.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'10.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.65%  "
$ws.Range("D26").Value = "'169.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "'20.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'5.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.93%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "'1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "'4.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").Value = "'0.0618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "'2.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("D37").Value = "'0.0838"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("D38").Value = "'1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").Value = "'5.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "'0.0955"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.96%  "
$ws.Range("D44").Value = "'97.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'16.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("D46").Value = "1.305.38"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "'2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.11%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "2.244.06"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'44.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.18%  "
